$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before row 868. This shifts the existing
# rows 868-936 down to 870-938 (dimension grows from A1:T936 to A1:T938),
# matching the diff exactly (verified: before-row N == after-row N+2 for
# every row from 868 to 936).
$ws.Rows("868:869").Insert()

# --- New row 868 ---
$ws.Range("A868").Value = 10
$ws.Range("B868").Value = "Vega Modelo de Temuco"
$ws.Range("C868").Value = "La Araucanía"
$ws.Range("D868").Value = 45106
$ws.Range("E868").Value = 9
$ws.Range("F868").Value = "Fruta"
$ws.Range("G868").Value = 100102
$ws.Range("H868").Value = "Cítricos"
$ws.Range("I868").Value = 100102004
$ws.Range("J868").Value = "Mandarina"
$ws.Range("K868").Value = "Clementina"
$ws.Range("L868").Value = "Primera"
$ws.Range("M868").Value = 8
$ws.Range("N868").Value = 315000
$ws.Range("O868").Value = 315000
$ws.Range("P868").Value = 315000
$ws.Range("Q868").Value = "$/bins (450 kilos)"
$ws.Range("R868").Value = "Región de O'Higgins"
$ws.Range("S868").Value = 700
$ws.Range("T868").Value = 450

# --- New row 869 ---
$ws.Range("A869").Value = 10
$ws.Range("B869").Value = "Vega Modelo de Temuco"
$ws.Range("C869").Value = "La Araucanía"
$ws.Range("D869").Value = 45106
$ws.Range("E869").Value = 9
$ws.Range("F869").Value = "Fruta"
$ws.Range("G869").Value = 100102
$ws.Range("H869").Value = "Cítricos"
$ws.Range("I869").Value = 100102004
$ws.Range("J869").Value = "Mandarina"
$ws.Range("K869").Value = "Clemenuless"
$ws.Range("L869").Value = "Primera"
$ws.Range("M869").Value = 200
$ws.Range("N869").Value = 10000
$ws.Range("O869").Value = 10000
$ws.Range("P869").Value = 10000
$ws.Range("Q869").Value = "$/bandeja 18 kilos"
$ws.Range("R869").Value = "Región de O'Higgins"
$ws.Range("S869").Value = 556
$ws.Range("T869").Value = 18
